$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 = "I0" and J1 = "IF" with the same style as H1 ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-8: I column = 1 (constant), J column = copy of H column value ---
foreach ($row in 2..8) {
    $hVal = $ws.Range("H$row").Value2
    $ws.Range("I$row").Value = 1
    $ws.Range("J$row").Value = $hVal
}
